# Update job title/company lines in PROFESSIONAL EXPERIENCE section
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "PRINCIPAL RESEARCH CONSULTANT - Clarity and Rigour, Washington, DC | 2012 – 2014"; New = "DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 – 2014" },
    @{ Old = "DIRECTOR OF RESEARCH - Helm, Washington, DC | 2010 – 2012"; New = "SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 – 2012" },
    @{ Old = "SENIOR RESEARCH ANALYST - GSD&M, Austin, TX | 2008 – 2010"; New = "SENIOR ANALYST - Myers Research, Washington, DC | 2008 – 2010" },
    @{ Old = "RESEARCH COORDINATOR - Salsa Labs, Inc., Washington, DC | 2004 – 2006"; New = "SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 – 2006" },
    @{ Old = "RESEARCH MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004"; New = "INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004" },
    @{ Old = "RESEARCH ANALYST - Lake Research Partners, Washington, DC | 2001 – 2002"; New = "PROGRAMMER - Lake Research Partners, Washington, DC | 2001 – 2002" },
    @{ Old = "FIELD RESEARCH COORDINATOR - The Feldman Group, Washington, DC | 2000 – 2001"; New = "FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 – 2001" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $found = $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $($r.Old)"
    }
}

$d.Save()
